# Rerun and summarise models without urban landuse:
#  - rename each summary sheet to its new "summNNNNNNNN" identifier
#  - relabel the "Education[T.Unknown]" coefficient row as
#    "Education[T.Unknown/Other]" on every sheet

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ30569822",
    "summ30800917",
    "summ31067981",
    "summ31332698",
    "summ31577454",
    "summ31866557",
    "summ32137466",
    "summ32498157",
    "summ32748135"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]

    if ($ws.Range("A5").Text -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
